$wb = $excel.ActiveWorkbook

# Row 65 data for each of the 4 worksheets (DE_LFT_#1, DE_LFT_#2, DE_PLT_#1, DE_PLT_#2)
# Columns: A=time B=总长 C=ID D=实际长度 E=和校验 F=总长_DEC G=ID_DEC H=实际长度_DEC I=和校验_DEC

$dateVal = 45851.43619212963

$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(65,1).Value = $dateVal
$ws1.Cells.Item(65,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(65,2).Value = "0x01,0x7c"
$ws1.Cells.Item(65,3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws1.Cells.Item(65,4).Value = "0x01,0x54"
$ws1.Cells.Item(65,5).Value = "0x14"
$ws1.Cells.Item(65,6).Value = 380
$ws1.Cells.Item(65,7).Value = [double]"7.598631275147109e+23"
$ws1.Cells.Item(65,8).Value = 340
$ws1.Cells.Item(65,9).Value = 14

$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(65,1).Value = $dateVal
$ws2.Cells.Item(65,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(65,2).Value = "0x01,0x7c"
$ws2.Cells.Item(65,3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws2.Cells.Item(65,4).Value = "0x01,0x54"
$ws2.Cells.Item(65,5).Value = "0xe"
$ws2.Cells.Item(65,6).Value = 380
$ws2.Cells.Item(65,7).Value = [double]"5.68432987514711e+23"
$ws2.Cells.Item(65,8).Value = 340
$ws2.Cells.Item(65,9).Value = 14

$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(65,1).Value = $dateVal
$ws3.Cells.Item(65,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(65,2).Value = "0x00,0x82"
$ws3.Cells.Item(65,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws3.Cells.Item(65,4).Value = "0x00,0x7C"
$ws3.Cells.Item(65,5).Value = "0x7"
$ws3.Cells.Item(65,6).Value = 130
$ws3.Cells.Item(65,7).Value = [double]"5.68631262647114e+23"
$ws3.Cells.Item(65,8).Value = 124
$ws3.Cells.Item(65,9).Value = 7

$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(65,1).Value = $dateVal
$ws4.Cells.Item(65,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws4.Cells.Item(65,2).Value = "0x00,0x82"
$ws4.Cells.Item(65,3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws4.Cells.Item(65,4).Value = "0x00,0x7B"
$ws4.Cells.Item(65,5).Value = "0x3"
$ws4.Cells.Item(65,6).Value = 130
$ws4.Cells.Item(65,7).Value = [double]"9.85046333984776e+23"
$ws4.Cells.Item(65,8).Value = 123
$ws4.Cells.Item(65,9).Value = 3
